# Updated cryptos list with latest prices and volume changes
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'27.158.39"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.08%  "

$ws.Range("D3").Value = "'1.638.05"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.47%  "

$ws.Range("E4").Value = "  +0.27%  "

$ws.Range("D5").Value = "'217.19"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.60%  "

$ws.Range("D6").Value = "'0.517"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.64%  "

$ws.Range("E7").Value = "  +0.28%  "

$ws.Range("E8").Value = "  -0.64%  "

$ws.Range("E9").Value = "  -0.09%  "

$ws.Range("D10").Value = "'20.10"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.98%  "

$ws.Range("D11").Value = "'0.0850"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.49%  "

$ws.Range("D12").Value = "'1.867.32"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.33%  "

$ws.Range("D13").Value = "'1.637.71"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.04%  "

$ws.Range("E14").Value = "  +0.16%  "

$ws.Range("E15").Value = "  +0.62%  "

$ws.Range("D16").Value = "'66.35"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.78%  "

$ws.Range("D17").Value = "'27.161.37"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.22%  "

$ws.Range("D18").Value = "'0.0₃0738"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.16%  "

$ws.Range("D19").Value = "'216.76"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.31%  "

$ws.Range("E20").Value = "  +0.21%  "

$ws.Range("E21").Value = "  +1.36%  "

$ws.Range("D22").Value = "'2.54"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +4.41%  "

$ws.Range("E23").Value = "  -0.55%  "

$ws.Range("D24").Value = "'9.14"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.12%  "

$ws.Range("D25").Value = "'147.75"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.21%  "

$ws.Range("E26").Value = "  +0.38%  "

$ws.Range("E27").Value = "  -0.22%  "

$ws.Range("E28").Value = "  -0.68%  "

$ws.Range("D29").Value = "'15.66"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.07%  "

$ws.Range("D30").Value = "'0.0508"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.16%  "

$ws.Range("E31").Value = "  -0.57%  "

$ws.Range("D32").Value = "'3.38"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.00%  "

$ws.Range("D33").Value = "'3.02"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.42%  "

$ws.Range("D34").Value = "'1.299.17"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.97%  "

$ws.Range("E36").Value = "  +1.08%  "

$ws.Range("E37").Value = "  -0.73%  "

$ws.Range("D38").Value = "'0.550"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.92%  "

$ws.Range("E39").Value = "  +1.00%  "

$ws.Range("E40").Value = "  +0.21%  "

$ws.Range("B41").Value = "MXToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D41").Value = "'2.28"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +5.69%  "

$ws.Range("B42").Value = "TrustWalletToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D42").Value = "'0.810"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.13%  "

$ws.Range("D43").Value = "'5.35"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.61%  "

$ws.Range("D44").Value = "'1.777.67"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.31%  "

$ws.Range("D45").Value = "'62.45"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.00%  "

$ws.Range("D46").Value = "'90.68"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.19%  "

$ws.Range("E47").Value = "  -1.32%  "

$ws.Range("D48").Value = "'0.0₆0108"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.53%  "

$ws.Range("E49").Value = "  +0.37%  "

$ws.Range("D50").Value = "'7.60"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.71%  "

$ws.Range("D51").Value = "'0.751"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +13.29%  "
